# Update worksheet with new TPM values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = 0.8203074518761176
$ws.Cells.Item(2, 10).Value = 0.8725723693674974
$ws.Cells.Item(2, 13).Value = 7.004922666666666
$ws.Cells.Item(2, 14).Value = 21.014768
$ws.Cells.Item(2, 15).Value = 0.134029393318039
$ws.Cells.Item(2, 16).Value = 0.1464771065395205
$ws.Cells.Item(2, 17).Value = 0.3823917234506666
$ws.Cells.Item(2, 18).Value = 3.441525511056
$ws.Cells.Item(2, 19).Value = 0.1099453101092225
$ws.Cells.Item(2, 20).Value = 0.1278118759112847
$ws.Cells.Item(3, 9).Value = 0.8203074518761176
$ws.Cells.Item(3, 10).Value = 0.8725723693674974
$ws.Cells.Item(3, 15).Value = 0.3796257919253833
$ws.Cells.Item(3, 16).Value = 0.4148827819958515
$ws.Cells.Item(3, 19).Value = 0.3114098660407644
$ws.Cells.Item(3, 20).Value = 0.3620152520958991
$ws.Cells.Item(4, 9).Value = 0.8203074518761176
$ws.Cells.Item(4, 10).Value = 0.8725723693674974
$ws.Cells.Item(4, 13).Value = 6.293636666666667
$ws.Cells.Item(4, 14).Value = 18.88091
$ws.Cells.Item(4, 15).Value = 0.12041993100245
$ws.Cells.Item(4, 16).Value = 0.1316036924905903
$ws.Cells.Item(4, 17).Value = 0.3435633319966667
$ws.Cells.Item(4, 18).Value = 3.09206998797
$ws.Cells.Item(4, 19).Value = 0.09878136675571771
$ws.Cells.Item(4, 20).Value = 0.1148337457740259
$ws.Cells.Item(5, 9).Value = 0.8203074518761176
$ws.Cells.Item(5, 10).Value = 0.8725723693674974
$ws.Cells.Item(5, 13).Value = 13.324299
$ws.Cells.Item(5, 14).Value = 26.648598
$ws.Cells.Item(5, 15).Value = 0.2549418168249328
$ws.Cells.Item(5, 16).Value = 0.1857460205306503
$ws.Cells.Item(5, 17).Value = 0.727360158111
$ws.Cells.Item(5, 18).Value = 4.364160948666
$ws.Cells.Item(5, 19).Value = 0.2091306721363286
$ws.Cells.Item(5, 20).Value = 0.1620768452350133
$ws.Cells.Item(6, 9).Value = 0.8203074518761176
$ws.Cells.Item(6, 10).Value = 0.8725723693674974
$ws.Cells.Item(6, 13).Value = 5.800427666666667
$ws.Cells.Item(6, 14).Value = 17.401283
$ws.Cells.Item(6, 15).Value = 0.110983066929195
$ws.Cells.Item(6, 16).Value = 0.1212903984433873
$ws.Cells.Item(6, 17).Value = 0.3166395458956667
$ws.Cells.Item(6, 18).Value = 2.849755913061
$ws.Cells.Item(6, 19).Value = 0.09104023683408455
$ws.Cells.Item(6, 20).Value = 0.1058346503512743
$ws.Cells.Item(7, 7).Value = 0.011958
$ws.Cells.Item(7, 8).Value = 0.023916
$ws.Cells.Item(7, 9).Value = 0.1796925481238824
$ws.Cells.Item(7, 10).Value = 0.1274276306325027
$ws.Cells.Item(7, 13).Value = 7.004922666666666
$ws.Cells.Item(7, 14).Value = 21.014768
$ws.Cells.Item(7, 15).Value = 0.134029393318039
$ws.Cells.Item(7, 16).Value = 0.1464771065395205
$ws.Cells.Item(7, 17).Value = 0.08376486524799999
$ws.Cells.Item(7, 18).Value = 0.502589191488
$ws.Cells.Item(7, 19).Value = 0.02408408320881647
$ws.Cells.Item(7, 20).Value = 0.01866523062823576
$ws.Cells.Item(8, 7).Value = 0.011958
$ws.Cells.Item(8, 8).Value = 0.023916
$ws.Cells.Item(8, 9).Value = 0.1796925481238824
$ws.Cells.Item(8, 10).Value = 0.1274276306325027
$ws.Cells.Item(8, 15).Value = 0.3796257919253833
$ws.Cells.Item(8, 16).Value = 0.4148827819958515
$ws.Cells.Item(8, 17).Value = 0.237256190736
$ws.Cells.Item(8, 18).Value = 1.423537144416
$ws.Cells.Item(8, 19).Value = 0.06821592588461889
$ws.Cells.Item(8, 20).Value = 0.0528675298999525
$ws.Cells.Item(9, 7).Value = 0.011958
$ws.Cells.Item(9, 8).Value = 0.023916
$ws.Cells.Item(9, 9).Value = 0.1796925481238824
$ws.Cells.Item(9, 10).Value = 0.1274276306325027
$ws.Cells.Item(9, 13).Value = 6.293636666666667
$ws.Cells.Item(9, 14).Value = 18.88091
$ws.Cells.Item(9, 15).Value = 0.12041993100245
$ws.Cells.Item(9, 16).Value = 0.1316036924905903
$ws.Cells.Item(9, 17).Value = 0.07525930726000001
$ws.Cells.Item(9, 18).Value = 0.45155584356
$ws.Cells.Item(9, 19).Value = 0.02163856424673235
$ws.Cells.Item(9, 20).Value = 0.01676994671656441
$ws.Cells.Item(10, 7).Value = 0.011958
$ws.Cells.Item(10, 8).Value = 0.023916
$ws.Cells.Item(10, 9).Value = 0.1796925481238824
$ws.Cells.Item(10, 10).Value = 0.1274276306325027
$ws.Cells.Item(10, 13).Value = 13.324299
$ws.Cells.Item(10, 14).Value = 26.648598
$ws.Cells.Item(10, 15).Value = 0.2549418168249328
$ws.Cells.Item(10, 16).Value = 0.1857460205306503
$ws.Cells.Item(10, 17).Value = 0.159331967442
$ws.Cells.Item(10, 18).Value = 0.637327869768
$ws.Cells.Item(10, 19).Value = 0.04581114468860425
$ws.Cells.Item(10, 20).Value = 0.02366917529563696
$ws.Cells.Item(11, 7).Value = 0.011958
$ws.Cells.Item(11, 8).Value = 0.023916
$ws.Cells.Item(11, 9).Value = 0.1796925481238824
$ws.Cells.Item(11, 10).Value = 0.1274276306325027
$ws.Cells.Item(11, 13).Value = 5.800427666666667
$ws.Cells.Item(11, 14).Value = 17.401283
$ws.Cells.Item(11, 15).Value = 0.110983066929195
$ws.Cells.Item(11, 16).Value = 0.1212903984433873
$ws.Cells.Item(11, 17).Value = 0.06936151403800001
$ws.Cells.Item(11, 18).Value = 0.416169084228
$ws.Cells.Item(11, 19).Value = 0.01994283009511043
$ws.Cells.Item(11, 20).Value = 0.01545574809211304
